# "integrated grid optimized hp" - remove the placeholder/zero battery entry (row 2)
# and shift the table up. The comment/link cells in column I (and the hyperlink
# itself) are re-attached to the rows they now belong to.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("bat")

# 1) Delete the obsolete all-zero placeholder row. Excel shifts every row
#    below it up by one and re-points the relative formulas automatically.
$ws.Rows.Item(2).Delete()

# 2) Remove the hyperlink that used to sit on (old) I2 - it will be re-created
#    on its new home (I7) below.
$ws.Hyperlinks.Delete()

# 3) Tidy up the "comment" column (I) so each note is attached to the row it
#    actually documents now that the table has shifted:
#      I2 (ex "109")  -> cleared, row no longer carries a note
#      I4 (ex "109")  -> "110" note
#      I5 (ex "110")  -> "111" note
#      I6 (ex "111")  -> cleared
#      I7             -> BYD storage link ("108"), restyled as a hyperlink
#      I8             -> installation markup note ("109")
$ws.Cells.Item(2,9).Value = $null

$ws.Cells.Item(4,9).Value = "https://pvspeicher.htw-berlin.de/wp-content/uploads/2017/03/WENIGER-2017_03-Vergleich-verschiedener-Kennzahlen-zur-Bewertung-der-energetischen-Performance-von-PV-Batteriesystemen.pdf"

$ws.Cells.Item(5,9).Value = "Wirkungsgraddiskussion"

$ws.Cells.Item(6,9).Value = $null

$ws.Cells.Item(7,9).Style = "Link"
$ws.Cells.Item(7,9).Value = "https://greenakku.de/Batterien/Lithium-Batterien/BYD-B-Box-H-9-0-fuer-SMA::1630.html"
$ws.Hyperlinks.Add($ws.Cells.Item(7,9), "https://greenakku.de/Batterien/Lithium-Batterien/BYD-B-Box-H-9-0-fuer-SMA::1630.html") | Out-Null

$ws.Cells.Item(8,9).Value = "Laut Mail von Batterienhersteller 20 Prozent Aufschlag auf Kaufpreis für Installation und Montage"
